$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out old data that is no longer part of the sheet
$ws.Cells.Clear()

# Set the new values
$ws.Range("A1").Value = "Sherin"
$ws.Range("A2").Value = "Dua"

# Match the selection state from the target file
$ws.Range("A2").Select()
